$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update text-like cell values (safe as plain strings; Excel will not
# reinterpret them as numbers because of multiple dots / percent signs / words).
$ws.Range("D2").Value = '98.168.58'
$ws.Range("E2").Value = '  +4.25%  '
$ws.Range("D3").Value = '3.354.04'
$ws.Range("E3").Value = '  +9.90%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("E5").Value = '  +9.90%  '
$ws.Range("E6").Value = '  +3.60%  '
$ws.Range("E7").Value = '  +8.88%  '
$ws.Range("E8").Value = '  +4.46%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").Value = '3.350.13'
$ws.Range("E10").Value = '  +9.87%  '
$ws.Range("E11").Value = '  +0.76%  '
$ws.Range("E12").Value = '  +2.68%  '
$ws.Range("D13").Value = '97.830.51'
$ws.Range("E13").Value = '  +4.32%  '
$ws.Range("E14").Value = '  +8.55%  '
$ws.Range("E15").Value = '  +4.64%  '
$ws.Range("D16").Value = '3.972.79'
$ws.Range("E16").Value = '  +9.57%  '
$ws.Range("E17").Value = '  +5.06%  '
$ws.Range("D18").Value = '3.353.31'
$ws.Range("E18").Value = '  +10.10%  '
$ws.Range("E19").Value = '  +4.80%  '
$ws.Range("E20").Value = '  +5.99%  '
$ws.Range("E21").Value = '  +12.27%  '
$ws.Range("E22").Value = '  +4.62%  '
$ws.Range("E23").Value = '  +12.61%  '
$ws.Range("E24").Value = '  +5.74%  '
$ws.Range("E25").Value = '  +4.98%  '
$ws.Range("E27").Value = '  +4.41%  '
$ws.Range("D28").Value = '3.531.36'
$ws.Range("E28").Value = '  +9.57%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("E30").Value = '  +6.02%  '
$ws.Range("E31").Value = '  +2.27%  '
$ws.Range("E32").Value = '  +3.18%  '
$ws.Range("E33").Value = '  -1.66%  '
$ws.Range("E34").Value = '  +4.47%  '
$ws.Range("E35").Value = '  +8.98%  '
$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("E36").Value = '  +13.73%  '
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("E38").Value = '  -2.31%  '
$ws.Range("E39").Value = '  +4.58%  '
$ws.Range("E40").Value = '  +3.57%  '
$ws.Range("E41").Value = '  +5.02%  '
$ws.Range("E42").Value = '  +3.27%  '
$ws.Range("E43").Value = '  -0.82%  '
$ws.Range("E44").Value = '  +7.86%  '
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("E46").Value = '  +19.69%  '
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("E48").Value = '  +8.49%  '
$ws.Range("E49").Value = '  +9.30%  '
$ws.Range("E50").Value = '  +4.04%  '
$ws.Range("E51").Value = '  +8.53%  '

# Update cells whose new text happens to look like a plain number
# (e.g. "1.18"). Force the cell to Text format first so Excel keeps
# the exact original string instead of converting it to a float, then
# clear the formatting again so no extra number-format style lingers
# on the cell (matching the original inline-string cells, which carry
# no style at all).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '255.40'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '625.26'
$ws.Range("D6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.18'
$ws.Range("D7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.388'
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.799'
$ws.Range("D11").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.95'
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000247'
$ws.Range("D15").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.49'
$ws.Range("D17").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.61'
$ws.Range("D19").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '485.11'
$ws.Range("D21").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000207'
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.14'
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.68'
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '87.97'
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.04'
$ws.Range("D27").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.123'
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.25'
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.21'
$ws.Range("D35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '517.48'
$ws.Range("D36").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.34'
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.95'
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '24.85'
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.449'
$ws.Range("D41").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.69'
$ws.Range("D43").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.781'
$ws.Range("D46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '161.20'
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.93'
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.37'
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.31'
$ws.Range("D50").ClearFormats()
